# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows
# on Sheet1 to reflect the corrected/re-annotated dialog acts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> (DAMSLTag, DialogAct)
$updates = @(
    @{ Row = 6;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 10; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 15; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 28; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 33; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 37; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 40; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 46; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 56; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 72; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 77; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 86; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 95; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
